# Checkout.xlsx update ("Code has been updated")
#
# Semantic edit (recovered from the OOXML diff):
#   - "Test Steps" becomes the active sheet/tab (was "Test Cases").
#   - "Test Cases" selection moves from D12 to B12 (and loses tabSelected).
#   - "Test Steps" selection moves from A17 to F12 (and gains tabSelected).
#   - On the "Test Steps" sheet, the wait-time value used by the two
#     `pause` steps (rows 11/13/26/28, column E) changes from "5" to "7".
#   - On the "Test Steps" sheet, the object name used by the two
#     `uniqueSelect` steps (rows 12/27, column E) changes from
#     "ProductDetails_Button" to "ProductItem_Link".
#
# Everything else (row 14-19 / 29-34 shared-string index churn in the raw
# XML) is just Excel's own shared-string table compacting around the two
# edits above -- not a separate content change.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Test Cases"
$ws2 = $wb.Worksheets.Item(2)   # "Test Steps"

# --- content edits -----------------------------------------------------
# Leading "'" forces these to stay text (shared-string) cells instead of
# being auto-coerced to numbers, matching the original t="s" cell type.
$ws2.Range("E12").Value = "'ProductItem_Link"
$ws2.Range("E27").Value = "'ProductItem_Link"

$ws2.Range("E11").Value = "'7"
$ws2.Range("E13").Value = "'7"
$ws2.Range("E26").Value = "'7"
$ws2.Range("E28").Value = "'7"

# --- selection / active-tab edits --------------------------------------
$ws1.Range("B12").Select()

$ws2.Activate()
$ws2.Range("F12").Select()
